$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '26.022.47'
Set-TextValue $ws.Range('E2') '  -2.39%  '
Set-TextValue $ws.Range('D3') '1.668.90'
Set-TextValue $ws.Range('E3') '  -1.58%  '
Set-TextValue $ws.Range('D4') '1.004'
Set-TextValue $ws.Range('E4') '  -0.13%  '
Set-TextValue $ws.Range('D5') '216.88'
Set-TextValue $ws.Range('E5') '  -1.33%  '
Set-TextValue $ws.Range('E6') '  -0.43%  '
Set-TextValue $ws.Range('E7') '  -0.14%  '
Set-TextValue $ws.Range('D8') '0.2654'
Set-TextValue $ws.Range('E8') '  +0.51%  '
Set-TextValue $ws.Range('D9') '0.06407'
Set-TextValue $ws.Range('E9') '  +1.91%  '
Set-TextValue $ws.Range('E10') '  -1.66%  '
Set-TextValue $ws.Range('D11') '0.07439'
Set-TextValue $ws.Range('E11') '  +1.21%  '
Set-TextValue $ws.Range('D12') '1.682.30'
Set-TextValue $ws.Range('E12') '  -1.03%  '
Set-TextValue $ws.Range('D13') '4.502'
Set-TextValue $ws.Range('E13') '  -0.34%  '
Set-TextValue $ws.Range('D14') '0.5846'
Set-TextValue $ws.Range('E14') '  +1.02%  '
Set-TextValue $ws.Range('D15') '0.000008562'
Set-TextValue $ws.Range('E15') '  +1.51%  '
Set-TextValue $ws.Range('D16') '64.48'
Set-TextValue $ws.Range('E16') '  -1.63%  '
Set-TextValue $ws.Range('D17') '26.072.24'
Set-TextValue $ws.Range('E17') '  -2.31%  '
Set-TextValue $ws.Range('D18') '4.948'
Set-TextValue $ws.Range('E18') '  -0.94%  '
Set-TextValue $ws.Range('E19') '  -0.10%  '
Set-TextValue $ws.Range('D20') '10.79'
Set-TextValue $ws.Range('E20') '  -1.91%  '
Set-TextValue $ws.Range('D21') '193.55'
Set-TextValue $ws.Range('E21') '  +3.77%  '
Set-TextValue $ws.Range('D22') '6.225'
Set-TextValue $ws.Range('E22') '  -0.31%  '
Set-TextValue $ws.Range('E23') '  -0.11%  '
Set-TextValue $ws.Range('D24') '145.02'
Set-TextValue $ws.Range('E24') '  +0.20%  '
Set-TextValue $ws.Range('D25') '7.615'
Set-TextValue $ws.Range('E25') '  +1.46%  '
Set-TextValue $ws.Range('D26') '0.1195'
Set-TextValue $ws.Range('E26') '  +3.41%  '
Set-TextValue $ws.Range('D27') '15.73'
Set-TextValue $ws.Range('E27') '  -0.55%  '
Set-TextValue $ws.Range('D28') '0.06431'
Set-TextValue $ws.Range('E28') '  +13.97%  '
Set-TextValue $ws.Range('E29') '  -1.57%  '
Set-TextValue $ws.Range('D30') '1.317'
Set-TextValue $ws.Range('E30') '  -1.44%  '
Set-TextValue $ws.Range('D31') '3.547'
Set-TextValue $ws.Range('E31') '  +1.25%  '
Set-TextValue $ws.Range('D32') '3.522'
Set-TextValue $ws.Range('E32') '  +0.86%  '
Set-TextValue $ws.Range('D33') '1.648'
Set-TextValue $ws.Range('E33') '  -0.08%  '
Set-TextValue $ws.Range('D34') '1.021'
Set-TextValue $ws.Range('E34') '  -0.02%  '
Set-TextValue $ws.Range('D35') '0.6109'
Set-TextValue $ws.Range('E35') '  +1.77%  '
Set-TextValue $ws.Range('E36') '  +0.41%  '
Set-TextValue $ws.Range('E37') '  +0.28%  '
Set-TextValue $ws.Range('D38') '6.255'
Set-TextValue $ws.Range('E38') '  +7.07%  '
Set-TextValue $ws.Range('D39') '0.01605'
Set-TextValue $ws.Range('E39') '  -0.52%  '
Set-TextValue $ws.Range('D40') '1.088.94'
Set-TextValue $ws.Range('E40') '  -1.25%  '
Set-TextValue $ws.Range('D41') '0.8608'
Set-TextValue $ws.Range('E42') '  +0.58%  '
Set-TextValue $ws.Range('D43') '100.68'
Set-TextValue $ws.Range('E43') '  +0.96%  '
Set-TextValue $ws.Range('D44') '1.816.13'
Set-TextValue $ws.Range('E44') '  -2.01%  '
Set-TextValue $ws.Range('E45') '  +1.27%  '
Set-TextValue $ws.Range('D46') '56.38'
Set-TextValue $ws.Range('E46') '  -0.31%  '
Set-TextValue $ws.Range('E47') '  +0.13%  '
Set-TextValue $ws.Range('D48') '8.062'
Set-TextValue $ws.Range('E48') '  -0.86%  '
Set-TextValue $ws.Range('D49') '0.05236'
Set-TextValue $ws.Range('E49') '  -0.08%  '
Set-TextValue $ws.Range('D50') '0.4286'
Set-TextValue $ws.Range('E50') '  -0.86%  '
Set-TextValue $ws.Range('D51') '6.025'
Set-TextValue $ws.Range('E51') '  +4.18%  '
